$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34 ("templates" / "need various template tests") becomes the first
# concrete template test: mark it EXISTS and rename the description.
$ws.Range("A34").Value = "√"
$ws.Range("F34").Value = "simple template"

# Insert a new row right after it for the second template test, pushing the
# remaining rows (icon / related / softBindingsMatched / soft binding) down
# by one.
$ws.Rows.Item(35).Insert()
$ws.Range("A35").Value = "√"
$ws.Range("B35").Value = "POSITIVE"
$ws.Range("C35").Value = "actions"
$ws.Range("E35").Value = "templates"
$ws.Range("F35").Value = "special ""all"" template"

# The "related" row (now shifted down to row 38) is also marked EXISTS.
$ws.Range("A38").Value = "√"

# Update the UI selection to match the author's saved state.
$ws.Range("A38").Select()
